$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC001")

# Move the value "run" from A3 up to A2, clearing A3.
$ws.Range("A2").Value = $ws.Range("A3").Value2
$ws.Range("A3").ClearContents()
